# Zeiterfassung_Jeckle_Lukas.xlsx - neue Zeitposten gebucht.
# Adds seven new time-tracking entries (rows 17-23) to the "Zeiterfassung" sheet.
#
# The rows are entered in the same chronological order the author actually
# typed them in (17, 18, 19, 21, 22, 23) and the "Coding" entry is inserted
# afterwards as a new row 20 - this reproduces the exact shared-string
# ordering seen in the target workbook (new strings 25-30 from the first
# pass, 31-32 from the later insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- first pass: rows 17-22 (chronological entry order) -------------------
$ws.Cells.Item(17, 1).Value2 = 45587
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = "Emergency-Meeting"
$ws.Cells.Item(17, 4).Value = "Emergency-Meeting abgehalten. (Mail von Management missverstanden)"

$ws.Cells.Item(18, 1).Value2 = 45587
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = "Planung"
$ws.Cells.Item(18, 4).Value = "Komplette Projekt-Planung von C++ auf Java umstrukturiert. (Mail von Management missverstanden)"

$ws.Cells.Item(19, 1).Value2 = 45588
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(19, 3).Value = "Planung"
$ws.Cells.Item(19, 4).Value = "Planung von Java wieder auf C++ umgestellt. (Mail von Management missverstanden)"

$ws.Cells.Item(20, 1).Value2 = 45588
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = "Planung"
$ws.Cells.Item(20, 4).Value = "Präsentation für Visions-Präsentation am 24.10.2024 mit Christian Kasper vorbereitet."

$ws.Cells.Item(21, 1).Value2 = 45589
$ws.Cells.Item(21, 2).Value = 1.5
$ws.Cells.Item(21, 3).Value = "Besprechung"
$ws.Cells.Item(21, 4).Value = "Teilnahme an den anderen Visions-Präsentationen und Präsentation der eigenen Visions-Präsentation."

$ws.Cells.Item(22, 1).Value2 = 45589
$ws.Cells.Item(22, 2).Value = 0.25
$ws.Cells.Item(22, 3).Value = "Planung"
$ws.Cells.Item(22, 4).Value = "Aufgaben-, Projekt-, und Zeitmanagementplanung."

# --- second pass: the "Coding" entry was booked later and inserted between
#     what are now rows 19 and 21, becoming the new row 20 ------------------
$ws.Rows.Item(20).Insert() | Out-Null

$ws.Cells.Item(20, 1).Value2 = 45588
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = "Coding"
$ws.Cells.Item(20, 4).Value = "C++ Projekt aufgesetzt und CodeQL Code-Scan des Projekts erstellt."

# --- formatting -------------------------------------------------------------
# Give the new date cells (column A) the same number format as the existing
# booking rows above them (numFmtId 14, "m/d/yyyy") by copying the format
# from row 16 rather than assigning NumberFormat directly (which would
# register a duplicate custom format).
$ws.Cells.Item(16, 1).Copy() | Out-Null
$ws.Range("A17:A23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Widen columns C and D so the new, longer booking texts fit.
$ws.Columns.Item(3).ColumnWidth = 26.166666666666668
$ws.Columns.Item(4).ColumnWidth = 95.41666666666667

# Leave the selection where the author ended up after typing the data.
$ws.Range("A20").Select() | Out-Null
